# wednesday titrations for blue tank and E5
# Append the new titration reading as row 72 on the CRMAccuracyData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 72

$ws.Cells.Item($row, 1).Value = 20220323
$ws.Cells.Item($row, 2).Value = 2227.4597991999999
$ws.Cells.Item($row, 3).Value = 2224.4699999999998
$ws.Cells.Item($row, 4).Formula = "=100*(B$row-C$row)/C$row"
$ws.Cells.Item($row, 5).Value = 180
$ws.Cells.Item($row, 6).Value = "CRM OPENED 20220318"

# Match the new selection left by Excel after the data entry (diff shows
# the sheetView selection move to G72, one cell past the new row's data).
[void]$ws.Range("G72").Select()
